$wb = $excel.ActiveWorkbook

# Source sheet that already carries the header/index cell style (style index 1:
# bold font, thin border all sides, center/top aligned) used consistently across
# all sheets in this workbook.
$srcWs = $wb.Worksheets.Item("discrepant roles")

# Create the new worksheet and move it to the end of the tab strip.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "attributions to add"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)
# Re-fetch by name: after Move(), the old $newSheet reference no longer tracks the sheet.
$ws = $wb.Worksheets.Item("attributions to add")

# Header row (row 1), starting at column B.
$headers = @("ID", "title", "role", "identification", "indicated_value", "text_bdrc_id", "text_84000_ids", "attribution_lang")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows (columns B..I); column A gets the running 0-based index.
$data = @(
    @("D51", "go cha’i bkod pa bstan pa", "translatorTib", "unknown", "Gö Chödrup", "WA0RK0051", "eft:g-ch-drup", "bo"),
    @("D84", "bu mo rnam dag dad pas zhus pa", "translatorTib", "unknown", "Gö Chödrup", "WA0RK0084", "eft:g-ch-drup", "bo"),
    @("D119", "yongs su mya ngan las 'das pa chen po'i mdo/", "translatorTib", "unknown", "wang phab zhwun (wang phan zhun)", "WA0RK0119", "eft:wang-phab-zhwun-wang-phan-zhun-", "bo"),
    @("D119", "yongs su mya ngan las 'das pa chen po'i mdo/", "translatorTib", "unknown", "dge ba'i blo gros", "WA0RK0119", "eft:dge-ba-i-blo-gros", "bo"),
    @("D119", "yongs su mya ngan las 'das pa chen po'i mdo/", "translatorTib", "unknown", "rgya mtsho'i sde", "WA0RK0119", "eft:rgya-mtsho-i-sde", "bo"),
    @("D267", "dpang skong phyag brgya pa", "translatorTib", "unknown", "Thönmi Sambhoṭa", "WA0RK0267", "eft:th-nmi-sambhota", "bo"),
    @("D287", "dam pa’i chos dran pa nye bar gzhag pa", "translatorTib", "unknown", "Tsultrim Gyaltsen", "WA0RK0287", "eft:tsultrim-gyaltsen", "bo"),
    @("D287", "dam pa’i chos dran pa nye bar gzhag pa", "translatorTib", "unknown", "Shang Buchikpa", "WA0RK0287", "eft:shang-buchikpa", "bo"),
    @("D287", "dam pa’i chos dran pa nye bar gzhag pa", "translatorTib", "unknown", "Sherap Ö", "WA0RK0287", "eft:sherap-", "bo"),
    @("D300", "dge ba’i bshes gnyen bsten pa’i mdo", "translatorTib", "unknown", "Paṇḍita Dharmākara", "WA0RK0300", "eft:pandita-dharmakara", "sa"),
    @("D300", "dge ba’i bshes gnyen bsten pa’i mdo", "translatorTib", "unknown", "Lotsāwa Zangkyong (bzang skyong)", "WA0RK0300", "eft:lotsawa-zangkyong-bzang-skyong-", "bo"),
    @("D312", "yangs pa’i grong khyer du ’jug pa’i mdo chen po", "translatorTib", "unknown", "Surendrabodhi", "WA0RK0312", "eft:surendrabodhi", "sa"),
    @("D312", "yangs pa’i grong khyer du ’jug pa’i mdo chen po", "translatorTib", "unknown", "Yeshé Dé", "WA0RK0312", "eft:yesh-d-", "bo"),
    @("D438", "sgrol ma la phyag ’tshal nyi shu rtsa gcig gis bstod pa", "translatorTib", "unknown", "Nyen Lotsawa Darma Drak", "WA0RK0438", "eft:nyen-lotsawa-darma-drak", "bo"),
    @("D674", "tshe dang ye shes dpag tu med pa’i mdo", "translatorTib", "unknown", "Patsap Nyima Drak [?]", "WA0RK0668", "eft:patsap-nyima-drak-", "bo"),
    @("D849", "tshe dang ye shes dpag tu med pa’i mdo", "translatorTib", "unknown", "Patsap Nyima Drak [?]", "WA0RK0668", "eft:patsap-nyima-drak-", "bo"),
    @("D830", "ye shes rngam pa glog gi 'khor lo/", "translatorTib", "unknown", "vajrvisramitra", "WA0RK0824", "eft:vajrvisramitra", "sa"),
    @("D830", "ye shes rngam pa glog gi 'khor lo/", "translatorTib", "unknown", "vairocanaraksita", "WA0RK0824", "eft:vairocanaraksita", "sa"),
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    $ws.Cells.Item($excelRow, 1).Value = $r
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 2).Value = $row[$c]
    }
}

$lastRow = $data.Length + 1

# Apply the shared "header/index" style to the header row (B1:I1) by copying
# format from an existing header cell, preserving the workbook-wide shared style
# index instead of minting a new one.
$srcWs.Range("B1").Copy()
$ws.Range("B1:I1").PasteSpecial(-4122)

# Apply the same style to column A (A1:A<lastRow>).
$srcWs.Range("A2").Copy()
$ws.Range("A1:A$lastRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0
